$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 40000
$ws.Range("D2").Value = 0.3624081495619552
$ws.Range("E2").Value = 2.830591467979062
$ws.Range("F2").Value = 0.5363
$ws.Range("H2").Value = 2.995599830918374

$ws.Range("B3").Value = 40000
$ws.Range("D3").Value = 0.5998284185563134
$ws.Range("E3").Value = 2.66221189799045
$ws.Range("F3").Value = 0.9236
$ws.Range("H3").Value = 2.995599830918374

$ws.Range("B4").Value = 40000
$ws.Range("D4").Value = 0.6677861195779499
$ws.Range("E4").Value = 2.623477043429551
$ws.Range("F4").Value = 1.142
$ws.Range("H4").Value = 2.995599830918374
